$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Text($cellRef, $text) {
    $c = $ws.Range($cellRef)
    $c.Value2 = "'" + $text
    $c.ClearFormats()
}

# Row 3
Set-Text "A3" "2025-05-28 21-11-13"
Set-Text "B3" "Daniel "
Set-Text "C3" "2000-01-01"
Set-Text "D3" "(917) 975-2625"
Set-Text "E3" "dshifrin5@gmail.com"
Set-Text "F3" "1181"
Set-Text "G3" "141"

# Row 4
Set-Text "A4" "2025-05-28 21-20-50"
Set-Text "B4" "Daniel Shifrin"
Set-Text "C4" "2002-04-25"
Set-Text "D4" "9179752625"
Set-Text "E4" "dshifrin@sandiego.edu"
Set-Text "F4" "1181"
Set-Text "G4" "Check up"

# Row 5
Set-Text "A5" "2025-05-28 21-37-01"
Set-Text "B5" "Daniel Shifrin"
Set-Text "C5" "2002-04-25"
Set-Text "D5" "(917) 975-2625"
Set-Text "E5" "dshifrin5@gmail.com"
Set-Text "F5" "1181"
Set-Text "G5" "Checkup"

# Row 6
Set-Text "A6" "2025-05-28 21-48-14"
Set-Text "B6" "Daniel "
Set-Text "C6" "2013-12-31"
Set-Text "D6" "9179752625"
Set-Text "E6" "dshifrin5@gmail.com"
Set-Text "F6" "1181"
Set-Text "G6" "ygtu"

# Row 7
Set-Text "A7" "2025-05-29 23-09-44"
Set-Text "E7" "robert99023@gmail.com"

# Row 8
Set-Text "A8" "2025-05-29 23-09-45"
Set-Text "E8" "robert99023@gmail.com"

# Row 9
Set-Text "A9" "2025-05-29 23-17-34"
Set-Text "E9" "elizavetakutko@gmail.com"

# Row 10
Set-Text "A10" "2025-05-29 23-17-35"
Set-Text "E10" "elizavetakutko@gmail.com"
